# Update the B4 values in each of the per-sheet report tables.
$wb = $excel.ActiveWorkbook

$updates = @{
    "SSBWE_v1_echoes_detection"        = 100
    "SSBWE_v1_mean_distance_error"     = 0.005111816370211811
    "SSBWE_v1_STD_distance_error"      = 0.005064691840527618
    "SSBWE_v1_mean_amplitude_error_1"  = -2.425717612959285
    "SSBWE_v1_STD_amplitude_error_1"   = 14.62152509368644
    "SSBWE_v1_mean_amplitude_error_2"  = -3.008500820176344
    "SSBWE_v1_STD_amplitude_error_2"   = 14.85535583866629
    "SSBWE_v2_echoes_detection"        = 100
    "SSBWE_v2_mean_distance_error"     = 0.005815927374713987
    "SSBWE_v2_STD_distance_error"      = 0.005149652958594186
    "SSBWE_v2_mean_amplitude_error_1"  = -6.620246088856059
    "SSBWE_v2_STD_amplitude_error_1"   = 13.4216449818085
    "SSBWE_v2_mean_amplitude_error_2"  = -7.244082523008549
    "SSBWE_v2_STD_amplitude_error_2"   = 13.62376233002742
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B4").Value = $updates[$sheetName]
}
